$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binary Search 2")

# Add new row 4: S.no=2, Page No. in notes="Binary S2 3",
# Question="Square root of N upto 3 decimal places", Link="na"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "Binary S2 3"
$ws.Range("D4").Value = "Square root of N upto 3 decimal places"
$ws.Range("E4").Value = "na"

# Match formatting of existing rows (style s="5" -> left/top aligned, font2)
$ws.Range("B4:E4").HorizontalAlignment = -4131
$ws.Range("B4:E4").VerticalAlignment = -4160

# Question cell needs wrap text to match the new cellXf (s="6")
$ws.Range("D4").WrapText = $true

$ws.Rows.Item(4).RowHeight = 28.8

$ws.Range("F4").Select()
